$wb = $excel.ActiveWorkbook

# Rename the existing "Sheet 1" to "Invoice"
$invoice = $wb.Worksheets.Item(1)
$invoice.Name = "Invoice"

# Add a new "Customer" worksheet after the Invoice sheet
$customer = $wb.Worksheets.Add($null, $invoice)
$customer.Name = "Customer"

# Populate the Customer sheet with header + data
$customer.Range("A1").Value = "customer_id"
$customer.Range("B1").Value = "customer_name"
$customer.Range("A2").Value = 2095
$customer.Range("B2").Value = "Marsha Billings"

# Make Customer the active sheet/tab
$customer.Activate()
